# Table B-1 (the first table in the document) has 3 columns:
#   1) stratum/category name (left-aligned, unchanged)
#   2) "Draft 2025 ADP" values (already right-aligned, unchanged)
#   3) "Final 2025 ADP" values (currently left-aligned -> change to right)
#
# Some rows are section-header rows whose single cell spans all three
# grid columns (gridSpan=3); those rows only have 1 cell and must be
# skipped since there is no "column 3" to adjust there.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

for ($r = 1; $r -le $t.Rows.Count; $r++) {
    $row = $t.Rows.Item($r)
    if ($row.Cells.Count -ge 3) {
        $cell = $t.Cell($r, 3)
        $cell.Range.ParagraphFormat.Alignment = 2  # wdAlignParagraphRight
    }
}
